$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 300
$ws.Range("I2").Value = 200
$ws.Range("J2").Value = 400
$ws.Range("K2").Value = 200
$ws.Range("L2").Value = 400
$ws.Range("M2").Value = -87
$ws.Range("N2").Value = -626
$ws.Range("H4").Value = 438
$ws.Range("I4").Value = 344
$ws.Range("J4").Value = 1002
$ws.Range("K4").Value = 344
$ws.Range("L4").Value = 1002
$ws.Range("M4").Value = -230
$ws.Range("N4").Value = -1230
$ws.Range("H9").Value = 176.66667
$ws.Range("I9").Value = 164.28572
$ws.Range("J9").Value = 194
$ws.Range("K9").Value = 164.28572
$ws.Range("L9").Value = 194
$ws.Range("M9").Value = 4.714280000000002
$ws.Range("N9").Value = -532
$ws.Range("H18").Value = 2125.8572
$ws.Range("I18").Value = 3322.75
$ws.Range("J18").Value = 530
$ws.Range("K18").Value = 3322.75
$ws.Range("L18").Value = 530
$ws.Range("M18").Value = -3038.75
$ws.Range("N18").Value = -1098
$ws.Range("H38").Value = 445.5
$ws.Range("J38").Value = 2000
$ws.Range("L38").Value = 6000
$ws.Range("N38").Value = -6744
$ws.Range("H39").Value = 350.33334
$ws.Range("I39").Value = 96.2
$ws.Range("K39").Value = 288.6
$ws.Range("M39").Value = 7.399999999999977
$ws.Range("H40").Value = 1880.9524
$ws.Range("I40").Value = 1600
$ws.Range("K40").Value = 1600
$ws.Range("M40").Value = -1425
$ws.Range("H51").Value = 4352.476
$ws.Range("J51").Value = 5557.143
$ws.Range("L51").Value = 5557.143
$ws.Range("N51").Value = -6525.143
$ws.Range("H52").Value = 0
$ws.Range("I52").Value = 0
$ws.Range("K52").Value = 0
$ws.Range("M52").ClearContents()
$ws.Range("H55").Value = 76
$ws.Range("I55").Value = 72
$ws.Range("J55").Value = 100
$ws.Range("K55").Value = 72
$ws.Range("L55").Value = 100
$ws.Range("M55").Value = 142
$ws.Range("N55").Value = -528
$ws.Range("H76").Value = 3122.7273
$ws.Range("I76").Value = 3114.2856
$ws.Range("K76").Value = 3114.2856
$ws.Range("M76").Value = -2799.2856
$ws.Range("H79").Value = 3122.7273
$ws.Range("I79").Value = 3114.2856
$ws.Range("K79").Value = 3114.2856
$ws.Range("M79").Value = -2022.2856
$ws.Range("H103").Value = 4232.5
$ws.Range("I103").Value = 6611.1113
$ws.Range("J103").Value = 1174.2858
$ws.Range("K103").Value = 19833.3339
$ws.Range("L103").Value = 3522.8574
$ws.Range("M103").Value = -19247.3339
$ws.Range("N103").Value = -4694.857400000001
$ws.Range("H132").Value = 1897.3835
$ws.Range("I132").Value = 792.20966
$ws.Range("J132").Value = 8126.5454
$ws.Range("K132").Value = 2376.62898
$ws.Range("L132").Value = 24379.6362
$ws.Range("M132").Value = 153.37102
$ws.Range("N132").Value = -29439.6362
$ws.Range("H137").Value = 2105.45
$ws.Range("I137").Value = 1396.0769
$ws.Range("J137").Value = 3422.8572
$ws.Range("K137").Value = 4188.2307
$ws.Range("L137").Value = 10268.5716
$ws.Range("M137").Value = -1638.2307
$ws.Range("N137").Value = -15368.5716

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2813.8462
$ws.Range("I132").Value = 2603.923
$ws.Range("J132").Value = 3233.6924
$ws.Range("K132").Value = 7811.768999999999
$ws.Range("L132").Value = 9701.0772
$ws.Range("M132").Value = -5281.768999999999
$ws.Range("N132").Value = -14761.0772

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H92").Value = 30000
$ws.Range("J92").Value = 30000
$ws.Range("L92").Value = 30000
$ws.Range("N92").Value = -34992
$ws.Range("H94").Value = 1528.8572
$ws.Range("I94").Value = 1578
$ws.Range("J94").Value = 1463.3334
$ws.Range("K94").Value = 1578
$ws.Range("L94").Value = 1463.3334
$ws.Range("M94").Value = -1127
$ws.Range("N94").Value = -2365.3334
$ws.Range("H95").Value = 12800
$ws.Range("J95").Value = 12800
$ws.Range("L95").Value = 12800
$ws.Range("N95").Value = -18292
$ws.Range("H96").Value = 3878
$ws.Range("J96").Value = 3878
$ws.Range("L96").Value = 3878
$ws.Range("N96").Value = -9370
$ws.Range("H99").Value = 3281.9443
$ws.Range("I99").Value = 2506.3333
$ws.Range("J99").Value = 4833.1665
$ws.Range("K99").Value = 2506.3333
$ws.Range("L99").Value = 4833.1665
$ws.Range("M99").Value = -1008.3333
$ws.Range("N99").Value = -7829.1665
$ws.Range("H126").Value = 3281.9443
$ws.Range("I126").Value = 2506.3333
$ws.Range("J126").Value = 4833.1665
$ws.Range("K126").Value = 7518.999899999999
$ws.Range("L126").Value = 14499.4995
$ws.Range("M126").Value = -5048.999899999999
$ws.Range("N126").Value = -19439.4995
$ws.Range("H132").Value = 4542.75
$ws.Range("I132").Value = 4340.375
$ws.Range("J132").Value = 4947.5
$ws.Range("K132").Value = 13021.125
$ws.Range("L132").Value = 14842.5
$ws.Range("M132").Value = -10491.125
$ws.Range("N132").Value = -19902.5
$ws.Range("H134").Value = 39287690
$ws.Range("I134").Value = 5002184.5
$ws.Range("J134").Value = 125001450
$ws.Range("K134").Value = 15006553.5
$ws.Range("L134").Value = 375004350
$ws.Range("M134").Value = -15004018.5
$ws.Range("N134").Value = -375009420

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 11459047
$ws.Range("I5").Value = 12500497
$ws.Range("J5").Value = 10417598
$ws.Range("K5").Value = 37501491
$ws.Range("L5").Value = 31252794
$ws.Range("M5").Value = -37501379
$ws.Range("N5").Value = -31253018
$ws.Range("H131").Value = 844.1
$ws.Range("I131").Value = 489
$ws.Range("J131").Value = 906.7646999999999
$ws.Range("K131").Value = 1467
$ws.Range("L131").Value = 2720.2941
$ws.Range("M131").Value = 3573
$ws.Range("N131").Value = -12800.2941
$ws.Range("H132").Value = 3630238.8
$ws.Range("I132").Value = 1668576.5
$ws.Range("J132").Value = 12348738
$ws.Range("K132").Value = 15017188.5
$ws.Range("L132").Value = 111138642
$ws.Range("M132").Value = -15014658.5
$ws.Range("N132").Value = -111143702
$ws.Range("H135").Value = 11459047
$ws.Range("I135").Value = 12500497
$ws.Range("J135").Value = 10417598
$ws.Range("K135").Value = 112504473
$ws.Range("L135").Value = 93758382
$ws.Range("M135").Value = -112501938
$ws.Range("N135").Value = -93763452

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H51").Value = 29250
$ws.Range("J51").Value = 29250
$ws.Range("L51").Value = 29250
$ws.Range("N51").Value = -30268

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3538.9546
$ws.Range("I7").Value = 3529.0527
$ws.Range("J7").Value = 3601.6667
$ws.Range("K7").Value = 3529.0527
$ws.Range("L7").Value = 3601.6667
$ws.Range("M7").Value = -3417.0527
$ws.Range("N7").Value = -3825.6667
$ws.Range("H22").Value = 40711.48
$ws.Range("I22").Value = 200356
$ws.Range("J22").Value = 800.35
$ws.Range("K22").Value = 200356
$ws.Range("L22").Value = 800.35
$ws.Range("M22").Value = -200061
$ws.Range("N22").Value = -1390.35
$ws.Range("H27").Value = 40711.48
$ws.Range("I27").Value = 200356
$ws.Range("J27").Value = 800.35
$ws.Range("K27").Value = 200356
$ws.Range("L27").Value = 800.35
$ws.Range("M27").Value = -200249
$ws.Range("N27").Value = -1014.35
$ws.Range("H108").Value = 28817
$ws.Range("J108").Value = 28817
$ws.Range("L108").Value = 28817
$ws.Range("N108").Value = -36497
$ws.Range("H123").Value = 25000
$ws.Range("J123").Value = 25000
$ws.Range("L123").Value = 25000
$ws.Range("N123").Value = -34800
$ws.Range("H126").Value = 3538.9546
$ws.Range("I126").Value = 3529.0527
$ws.Range("J126").Value = 3601.6667
$ws.Range("K126").Value = 10587.1581
$ws.Range("L126").Value = 10805.0001
$ws.Range("M126").Value = -8117.158100000001
$ws.Range("N126").Value = -15745.0001
$ws.Range("H127").Value = 35000
$ws.Range("J127").Value = 35000
$ws.Range("L127").Value = 35000
$ws.Range("N127").Value = -44920

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H128").Value = 31458.7
$ws.Range("J128").Value = 31458.7
$ws.Range("L128").Value = 31458.7
$ws.Range("N128").Value = -41418.7
$ws.Range("H131").Value = 20000
$ws.Range("J131").Value = 20000
$ws.Range("L131").Value = 20000
$ws.Range("N131").Value = -30080
